# Appends rows 45-57 (new car listings) to the worksheet, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('Dacia', 'Sandero', '€ 1 000,-5', 'Pas d''information', '100 326 km', 'Boîte manuelle', '10/2012', 'Essence', '55 kW (75 CH)', 'Professionnel', '5', '5', '55 kW (75 CH) Boîte manuelle 1 149 cm³ 5 4 1 050 kg ', 'https://www.autoscout24.be/fr/offres/dacia-sandero-1-2i-75pk-schadewagen-airco-blueto-essence-gris-5a4524ae-2307-499f-bc27-9601fba15f4e'),
    @('Peugeot', '208', '€ 5 990,-5', 'Très bon prix', '118 000 km', 'Boîte manuelle', '01/2016', 'Diesel', '55 kW (75 CH)', 'Professionnel', '5', '5', '55 kW (75 CH) Boîte manuelle 1 560 cm³ 5 4 1 140 kg ', 'https://www.autoscout24.be/fr/offres/peugeot-208-1-6-bluehdi-style-diesel-blanc-bd52d51c-57b3-4c08-b2ee-b0dcc25c5be4'),
    @('Volkswagen', 'Polo', '€ 6 990,-5', 'Très bon prix', '75 000 km', 'Boîte manuelle', '03/2013', 'Essence', '51 kW (69 CH)', 'Professionnel', '5', '5', '51 kW (69 CH) Boîte manuelle 1 198 cm³ 5 3 1 067 kg ', 'https://www.autoscout24.be/fr/offres/volkswagen-polo-1-2-essence-70cv-climatisation-essence-blanc-7419a210-815e-4abd-8263-5c1cbe07091e'),
    @('Opel', 'Insignia', '€ 7 250,-5', 'Très bon prix', '128 548 km', 'Boîte manuelle', '03/2016', 'Diesel', '100 kW (136 CH)', 'Professionnel', '5', '5', '100 kW (136 CH) Boîte manuelle 1 598 cm³ 6 4 1 613 kg ', 'https://www.autoscout24.be/fr/offres/opel-insignia-sports-tourer-1-6-cdti-ecoflex-cosmo-diesel-blanc-f8b601d1-1bcf-4bb8-b08d-3ce6351d1858'),
    @('Peugeot', 'Partner', '€ 1 950,-5', 'Pas d''information', '131 000 km', 'Boîte manuelle', '06/2004', 'Essence', '80 kW (109 CH)', 'Professionnel', '5', '5', '80 kW (109 CH) Boîte manuelle 1 587 cm³ 5 4 1 260 kg ', 'https://www.autoscout24.be/fr/offres/peugeot-partner-1-6i-benzine-bj-2004-131000km-5-zitplaatsen-airco-essence-rouge-66044004-8b6f-4b72-aace-b74bfb15e962'),
    @('Citroen', 'C3', '€ 10 250,-5', 'Très bon prix', '19 599 km', 'Boîte manuelle', '10/2022', 'Essence', '61 kW (83 CH)', 'Professionnel', '5', '5', '61 kW (83 CH) Boîte manuelle 1 199 cm³ 5 3 980 kg ', 'https://www.autoscout24.be/fr/offres/citroen-c3-1-2i-puretech-you-s-essence-rouge-7a6c2a29-4e6f-4016-b6e6-162b32206bc4'),
    @('Mercedes-Benz', 'B 180', '€ 2 750,-5', 'Prix correct', '189 000 km', 'Boîte automatique', '06/2005', 'Diesel', '80 kW (109 CH)', 'Professionnel', '5', '5', '80 kW (109 CH) Boîte automatique 1 991 cm³ 1 4 1 470 kg ', 'https://www.autoscout24.be/fr/offres/mercedes-benz-b-180-cdi-automatique-bj-2005-189000km-diesel-brun-08a6a9c1-263f-430e-87d3-e525b30c7785'),
    @('Opel', 'Insignia', '€ 3 600,-5', 'Très bon prix', '200 240 km', 'Boîte manuelle', '03/2013', 'Diesel', '96 kW (131 CH)', 'Professionnel', '5', '5', '96 kW (131 CH) Boîte manuelle 1 956 cm³ 6 4 1 610 kg ', 'https://www.autoscout24.be/fr/offres/opel-insignia-sports-tourer-2-0-cdti-cosmo-dpf-diesel-gris-6ad70203-4e1e-4928-94a4-2039f729ad3e'),
    @('Peugeot', '208', '€ 4 950,-5', 'Très bon prix', '118 431 km', 'Boîte manuelle', '04/2018', 'Essence', '50 kW (68 CH)', 'Professionnel', '5', '5', '50 kW (68 CH) Boîte manuelle 1 199 cm³ 5 3 1 035 kg ', 'https://www.autoscout24.be/fr/offres/peugeot-208-1-2i-puretech-access-essence-rouge-ee538c2c-2a94-4294-8e16-015a1e66a544'),
    @('Skoda', 'Rapid/Spaceback', '€ 7 450,-5', 'Très bon prix', '100 063 km', 'Boîte manuelle', '08/2015', 'Diesel', '77 kW (105 CH)', 'Professionnel', '5', '5', '77 kW (105 CH) Boîte manuelle 1 598 cm³ 5 4 1 265 kg ', 'https://www.autoscout24.be/fr/offres/skoda-rapid-spaceback-rapid-1-6-cr-tdi-elegance-diesel-noir-b81f55d8-98c2-46c7-98d7-e1bd9d29068e'),
    @('Mazda', '2', '€ 2 950,-5', 'Très bon prix', '149 218 km', 'Boîte manuelle', '06/2007', 'Essence', '55 kW (75 CH)', 'Professionnel', '5', '5', '55 kW (75 CH) Boîte manuelle 1 242 cm³ 5 4 1 125 kg ', 'https://www.autoscout24.be/fr/offres/mazda-2-1-25-essence-75cv-essence-argent-392f5a5e-0038-4f56-bdb8-b48b96e58351'),
    @('BMW', '216', '€ 9 950,-5', 'Très bon prix', '138 000 km', 'Boîte manuelle', '01/2018', 'Essence', '75 kW (102 CH)', 'Professionnel', '5', '5', '75 kW (102 CH) Boîte manuelle 1 499 cm³ 6 3 1 440 kg ', 'https://www.autoscout24.be/fr/offres/bmw-216-i-active-tourer-euro-6b-gps-navi-bluetooth-essence-gris-24ae4f29-68b4-4b32-97a1-7ce14aba3332'),
    @('Peugeot', '308', '€ 5 990,-5', 'Très bon prix', '118 000 km', 'Boîte manuelle', '01/2011', 'Essence', '88 kW (120 CH)', 'Professionnel', '4', '2', '88 kW (120 CH) Boîte manuelle 1 598 cm³ 5 4 1 555 kg ', 'https://www.autoscout24.be/fr/offres/peugeot-308-cc-1-6i-sport-pack-airco-garantie-prete-immatricul-essence-e57bcf32-9adc-4b18-a130-92b738dc1a62'),
)

$startRow = 45
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $val = $rowValues[$c]
        if ($val -match '^-?\d+(\.\d+)?$') {
            # Force pure-numeric-looking text (e.g. "5") to stay text, like the source inlineStr cells
            $cell.NumberFormat = "@"
            $cell.Value = $val
        } else {
            $cell.Value = $val
        }
    }
}
